$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add the new contact row data (rows 6-12)
$ws.Range("A6").Value = "mercadolibre4"
$ws.Range("B6").Value = "ml.com"
$ws.Range("C6").Value = "compras online"

$ws.Range("A7").Value = "MC6"
$ws.Range("A8").Value = "MC7"
$ws.Range("A9").Value = "jhjhjhhj"
$ws.Range("A10").Value = "ghfghfg"
$ws.Range("A11").Value = "thrtyrtyrt"

$ws.Range("A12").Value = "BUKIS 4444"
$ws.Range("D12").Value = "MARCO ANTONIO 444"

# Update the active selection to D12 as in the edited file
$ws.Range("D12").Select() | Out-Null
